$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")
$ws.Activate()

# Row 13 (Iteration 12) - Grant Type Profile
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "Enter Text to delete Grant Type Profile"
$ws.Range("C13").Value = "NewTitleEntered"

# Row 14 (Iteration 13) - Grant Category Profile
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Enter Text to delete Grant Category Profile"
$ws.Range("C14").Value = "NewTitleEntered"

# Match C13/C14 style with C9/C10 (NewTitleEntered cells) which use a distinct style
$ws.Range("C9").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active selection to match the post-edit state
$ws.Range("D15").Select()
